$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13, pushing existing rows 13.. down by one.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new data record.
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value = 45222
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = "Espárragos"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 1900
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 1950
$ws.Range("N13").Value = '$/kilo'
$ws.Range("O13").Value = "Provincia de Linares"
$ws.Range("P13").Value = 1950
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
